$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (first Iceland bank record) ---
$ws.Range("D2").Value = -0.0745
$ws.Range("E2").Value = -0.34
$ws.Range("F2").Value = 0.0638
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 28.4
$ws.Range("L2").Value = 0.08964646464646464
$ws.Range("M2").Value = 55.4
$ws.Range("N2").Value = 0.04315650074004829
$ws.Range("O2").Value = 1.950704225352113
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 55.4
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 834.4
$ws.Range("V2").Value = 0.6499961050089584
$ws.Range("W2").Value = 0.01792929292929293
$ws.Range("X2").Value = 0.07889947230462449
$ws.Range("Y2").Value = -0.06097017937533156
$ws.Range("Z2").Value = 0.08023320307256698
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04422734009395678
$ws.Range("AC2").Value = -0.04422734009395678
$ws.Range("AD2").Value = 2600.3
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2600.3
$ws.Range("AG2").Value = 1765.9
$ws.Range("AH2").Value = 0.6694902162718847
$ws.Range("AI2").Value = 0.6515899466259052
$ws.Range("AJ2").Value = 0.5790595487932844
$ws.Range("AK2").Value = 0.559484206190793

# AN2 / AP2 are removed entirely in the target workbook
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 (second Iceland bank record) ---
$ws.Range("B3").Value = "Arion banki hf. (ICSE:ARION)"

$ws.Range("D3").Value = -0.0745
$ws.Range("E3").Value = -0.34
$ws.Range("F3").Value = 0.0638
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 28.4
$ws.Range("L3").Value = 0.08964646464646464
$ws.Range("M3").Value = 55.4
$ws.Range("N3").Value = 0.04315650074004829
$ws.Range("O3").Value = 1.950704225352113
$ws.Range("P3").Value = -0.0
$ws.Range("Q3").Value = -0.0
$ws.Range("R3").Value = -0.0
$ws.Range("S3").Value = 55.4
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 834.4
$ws.Range("V3").Value = 0.6499961050089584
$ws.Range("W3").Value = 0.01792929292929293
$ws.Range("X3").Value = 0.07889947230462449
$ws.Range("Y3").Value = -0.06097017937533156
$ws.Range("Z3").Value = 0.08023320307256698
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04422734009395678
$ws.Range("AC3").Value = -0.04422734009395678
$ws.Range("AD3").Value = 2600.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 2600.3
$ws.Range("AG3").Value = 1765.9
$ws.Range("AH3").Value = 0.6694902162718847
$ws.Range("AI3").Value = 0.6515899466259052
$ws.Range("AJ3").Value = 0.5790595487932844
$ws.Range("AK3").Value = 0.559484206190793

# AN3 / AP3 are removed entirely in the target workbook
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
